{"js": "// Replace the 25 \"two-digit \u00f7 one-digit\" answer strings in the worksheet\n// table with their updated values, in document order. The table has 20\n// rows x 5 columns; only rows 0, 4, 8, 12, 16 carry text (the intervening\n// rows are blank spacer rows), so we walk those five rows left-to-right.\nconst replacements = [\n  \"26\u00f73=8, 2\", \"77\u00f73=25, 2\", \"35\u00f79=3, 8\", \"33\u00f79=3, 6\", \"61\u00f72=30, 1\",\n  \"42\u00f76=7, 0\", \"13\u00f72=6, 1\", \"26\u00f79=2, 8\", \"65\u00f72=32, 1\", \"73\u00f79=8, 1\",\n  \"49\u00f77=7, 0\", \"53\u00f74=13, 1\", \"47\u00f75=9, 2\", \"14\u00f75=2, 4\", \"49\u00f72=24, 1\",\n  \"99\u00f72=49, 1\", \"29\u00f78=3, 5\", \"42\u00f75=8, 2\", \"98\u00f77=14, 0\", \"28\u00f79=3, 1\",\n  \"59\u00f77=8, 3\", \"87\u00f78=10, 7\", \"52\u00f78=6, 4\", \"11\u00f74=2, 3\", \"81\u00f79=9, 0\",\n];\n\nconst dataRows = [0, 4, 8, 12, 16];\nconst cols = 5;\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nlet k = 0;\nfor (const r of dataRows) {\n  for (let c = 0; c < cols; c++) {\n    table.getCell(r, c).value = replacements[k];\n    k++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the 25 \"two-digit \u00f7 one-digit\" answer strings in the worksheet\n# table with their updated values, in document order. The table has 20\n# rows x 5 columns; only rows 1, 5, 9, 13, 17 (1-based) carry text (the\n# intervening rows are blank spacer rows), so we walk those five rows\n# left-to-right, same order as the source diff.\n\n$replacements = @(\n  \"26\u00f73=8, 2\", \"77\u00f73=25, 2\", \"35\u00f79=3, 8\", \"33\u00f79=3, 6\", \"61\u00f72=30, 1\",\n  \"42\u00f76=7, 0\", \"13\u00f72=6, 1\", \"26\u00f79=2, 8\", \"65\u00f72=32, 1\", \"73\u00f79=8, 1\",\n  \"49\u00f77=7, 0\", \"53\u00f74=13, 1\", \"47\u00f75=9, 2\", \"14\u00f75=2, 4\", \"49\u00f72=24, 1\",\n  \"99\u00f72=49, 1\", \"29\u00f78=3, 5\", \"42\u00f75=8, 2\", \"98\u00f77=14, 0\", \"28\u00f79=3, 1\",\n  \"59\u00f77=8, 3\", \"87\u00f78=10, 7\", \"52\u00f78=6, 4\", \"11\u00f74=2, 3\", \"81\u00f79=9, 0\"\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$dataRows = @(1, 5, 9, 13, 17)\n$k = 0\nforeach ($r in $dataRows) {\n  for ($c = 1; $c -le 5; $c++) {\n    $t.Cell($r, $c).Range.Text = $replacements[$k]\n    $k++\n  }\n}\n"}
